$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), matching the style of the
# existing header cells (bold, centered, bordered - same as H1).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New data columns I (I0) and J (IF) for rows 2-18.
$data = @(
    @(2, 9, 9),
    @(3, 8, 8),
    @(4, 3, 4),
    @(5, 8, 9),
    @(6, 9, 9),
    @(7, 5, 5),
    @(8, 10, 10),
    @(9, 8, 8),
    @(10, 5, 7),
    @(11, 8, 9),
    @(12, 7, 8),
    @(13, 6, 8),
    @(14, 8, 9),
    @(15, 8, 8),
    @(16, 1, 1),
    @(17, 3, 4),
    @(18, 1, 1)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $iVal = $entry[1]
    $jVal = $entry[2]
    $ws.Cells.Item($row, 9).Value = $iVal
    $ws.Cells.Item($row, 10).Value = $jVal
}
